$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 8682
    $ws.Range("F4").Value = 406
    $ws.Range("F5").Value = 32
}
